# The "invalid_logins" test-data sheet contained some rows of bogus login
# data that a parallel test run was picking up and failing on. Switch to
# that sheet, select the cell the author ended up on, drop the
# hyperlinks pointing into the rows we're about to remove, then delete
# rows 2-5 (the "ft"/roshan@ur.co.nz, fttt/232@3746 test rows) so the
# remaining rows shift up.
$wb  = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("invalid_logins")

$ws2.Activate()
$ws2.Range("D9").Select()

# Drop the hyperlinks that live in the rows being removed before the
# delete so no dangling hyperlink relationships are left behind.
$ws2.Range("A3").Hyperlinks.Delete()

# Remove the four rows of now-unwanted test data; everything below
# shifts up to close the gap.
$ws2.Rows("2:5").Delete()
